$d = $word.ActiveDocument

function Set-ParagraphRuns($paraIndex, $innerXml) {
    $para = $d.Paragraphs($paraIndex)
    $r = $para.Range
    # Exclude the trailing paragraph mark so w:pPr / paragraph identity is untouched.
    $target = $d.Range($r.Start, $r.End - 1)
    $pkg = '<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?>' + `
        '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' + `
        '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' + `
        '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p>' + `
        $innerXml + `
        '</w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $target.InsertXML($pkg)
}

# 1. "User: Admin" -> "User: Admin can"
$xml1 = '<w:r><w:t>User: Admin</w:t></w:r><w:r><w:t xml:space="preserve"> can</w:t></w:r>'
Set-ParagraphRuns 3 $xml1

# 2. "User: Bidder" -> "User: Bidder can"
$xml2 = '<w:r><w:t>User: Bidder</w:t></w:r><w:r><w:t xml:space="preserve"> can</w:t></w:r>'
Set-ParagraphRuns 12 $xml2

# 3. "300 concurrent users" -> "Site can 300 concurrent users"
$xml3 = '<w:r><w:t xml:space="preserve">Site can </w:t></w:r><w:r><w:t>300 concurrent users</w:t></w:r>'
Set-ParagraphRuns 19 $xml3

# 4. "Align project with Material Design specification" -> "Site should be aligned with Material Design specification"
#    (split into several runs, with the _GoBack bookmark re-anchored here)
$xml4 = '<w:r><w:t>Site should be a</w:t></w:r>' + `
    '<w:r><w:t>lign</w:t></w:r>' + `
    '<w:r><w:t>ed</w:t></w:r>' + `
    '<w:r><w:t xml:space="preserve"> </w:t></w:r>' + `
    '<w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/>' + `
    '<w:r><w:t>with Material Design specification</w:t></w:r>'
Set-ParagraphRuns 20 $xml4

# 5. Remove the old _GoBack bookmark that used to sit before " bid items"
#    (it moved to the "aligned" paragraph above).
$d.Bookmarks("_GoBack").Delete()
